$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.017.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.849.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.014"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'1.012"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'309.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4783"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3674"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.32%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.07232"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.9312"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'19.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.07741"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.843.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'5.347"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'6.457"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'89.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'1.016"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.000008667"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'1.012"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'27.031.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'14.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.071"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'10.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'153.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "  +1.96%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'2.015"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'114.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'4.966"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'3.303"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.54%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'1.179"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'0.7405"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'4.505"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'2.767"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.82%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'1.114"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.58%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.01959"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.04%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.05268"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.44%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'2.978"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.19%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.5263"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.34%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'7.014"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.00%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "  +5.92%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'8.244"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.4757"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.78%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'1.012"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'101.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.66%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.612"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.50%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'65.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.28%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.06077"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.8887"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.02%  "
$ws.Range("E51").Style = "Normal"
